$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 437.25
$ws.Range("I2").Value = 83
$ws.Range("K2").Value = 83
$ws.Range("M2").Value = 30

$ws.Range("H6").Value = 637624.75
$ws.Range("I6").Value = 637624.75
$ws.Range("K6").Value = 1912874.25
$ws.Range("M6").Value = -1912762.25

$ws.Range("H12").Value = 117.625
$ws.Range("I12").Value = 113.8
$ws.Range("J12").Value = 124
$ws.Range("K12").Value = 113.8
$ws.Range("L12").Value = 124
$ws.Range("M12").Value = 56.2
$ws.Range("N12").Value = -464

$ws.Range("H29").Value = 926.5
$ws.Range("I29").Value = 1503
$ws.Range("J29").Value = 350
$ws.Range("K29").Value = 4509
$ws.Range("L29").Value = 1050
$ws.Range("M29").Value = -4228
$ws.Range("N29").Value = -1612

$ws.Range("H38").Value = 2800.1428
$ws.Range("I38").Value = 25.25
$ws.Range("J38").Value = 6500
$ws.Range("K38").Value = 75.75
$ws.Range("L38").Value = 19500
$ws.Range("M38").Value = 296.25
$ws.Range("N38").Value = -20244

$ws.Range("H58").Value = 689.25
$ws.Range("J58").Value = 1043.8
$ws.Range("L58").Value = 3131.4
$ws.Range("N58").Value = -3431.4

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H125").Value = 1497.5
$ws.Range("I125").Value = 1495
$ws.Range("K125").Value = 13455
$ws.Range("M125").Value = -10995

$ws.Range("H135").Value = 2865
$ws.Range("I135").Value = 2350
$ws.Range("J135").Value = 3895
$ws.Range("K135").Value = 21150
$ws.Range("L135").Value = 35055
$ws.Range("M135").Value = -18615
$ws.Range("N135").Value = -40125

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 39
$ws.Range("I5").Value = 41.25
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 41.25
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 70.75
$ws.Range("N5").Value = -254

$ws.Range("H30").Value = 7497.5
$ws.Range("I30").Value = 9995
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 9995
$ws.Range("L30").Value = 5000
$ws.Range("M30").Value = -9845
$ws.Range("N30").Value = -5300

$ws.Range("H88").Value = 1350.4445
$ws.Range("I88").Value = 1109
$ws.Range("J88").Value = 1833.3334
$ws.Range("K88").Value = 1109
$ws.Range("L88").Value = 1833.3334
$ws.Range("M88").Value = -703
$ws.Range("N88").Value = -2645.3334

$ws.Range("H91").Value = 1350.4445
$ws.Range("I91").Value = 1109
$ws.Range("J91").Value = 1833.3334
$ws.Range("K91").Value = 1109
$ws.Range("L91").Value = 1833.3334
$ws.Range("M91").Value = 295
$ws.Range("N91").Value = -4641.3334

$ws.Range("H97").Value = 1408.2858
$ws.Range("I97").Value = 1408.2858
$ws.Range("K97").Value = 1408.2858
$ws.Range("M97").Value = -912.2858000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 39
$ws.Range("I4").Value = 41.25
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 41.25
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 73.75
$ws.Range("N4").Value = -260

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H94").Value = 2250
$ws.Range("I94").Value = 2250
$ws.Range("K94").Value = 2250
$ws.Range("M94").Value = -1799

$ws.Range("H99").Value = 1274.875
$ws.Range("I99").Value = 1299.8572
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 1299.8572
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 198.1428000000001
$ws.Range("N99").Value = -4096

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 165.83333
$ws.Range("I5").Value = 100.5
$ws.Range("K5").Value = 100.5
$ws.Range("M5").Value = 11.5

$ws.Range("H22").Value = 299
$ws.Range("I22").Value = 221.42857
$ws.Range("J22").Value = 480
$ws.Range("K22").Value = 221.42857
$ws.Range("L22").Value = 480
$ws.Range("M22").Value = 128.57143
$ws.Range("N22").Value = -1180

$ws.Range("H31").Value = 4823.2593
$ws.Range("I31").Value = 3291.3333
$ws.Range("J31").Value = 7887.1113
$ws.Range("K31").Value = 3291.3333
$ws.Range("L31").Value = 7887.1113
$ws.Range("M31").Value = -2996.3333
$ws.Range("N31").Value = -8477.1113

$ws.Range("H34").Value = 4823.2593
$ws.Range("I34").Value = 3291.3333
$ws.Range("J34").Value = 7887.1113
$ws.Range("K34").Value = 3291.3333
$ws.Range("L34").Value = 7887.1113
$ws.Range("M34").Value = -3089.3333
$ws.Range("N34").Value = -8291.1113

$ws.Range("H99").Value = 4500
$ws.Range("I99").Value = 4500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3002
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11030
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 874.75
$ws.Range("I13").Value = 374.5
$ws.Range("K13").Value = 1123.5
$ws.Range("M13").Value = -955.5

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H34").Value = 1924.8334
$ws.Range("I34").Value = 25
$ws.Range("J34").Value = 2874.75
$ws.Range("K34").Value = 75
$ws.Range("L34").Value = 8624.25
$ws.Range("M34").Value = 9
$ws.Range("N34").Value = -8792.25

$ws.Range("H39").Value = 3695
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H117").Value = 2185.5557
$ws.Range("I117").Value = 1939.6666
$ws.Range("K117").Value = 5818.9998
$ws.Range("M117").Value = -2376.9998

$ws.Range("H122").Value = 978
$ws.Range("J122").Value = 1019.6
$ws.Range("L122").Value = 9176.4
$ws.Range("N122").Value = -14076.4

$ws.Range("H129").Value = 1444.1666
$ws.Range("I129").Value = 332.5
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 997.5
$ws.Range("L129").Value = 6000
$ws.Range("M129").Value = 4002.5
$ws.Range("N129").Value = -16000

$ws.Range("H132").Value = 1490
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2199.8
$ws.Range("I80").Value = 1666.6666
$ws.Range("K80").Value = 1666.6666
$ws.Range("M80").Value = -668.6666

$ws.Range("H83").Value = 2199.8
$ws.Range("I83").Value = 1666.6666
$ws.Range("K83").Value = 8333.333000000001
$ws.Range("M83").Value = -3341.333000000001

$ws.Range("H97").Value = 549.5
$ws.Range("I97").Value = 549.5
$ws.Range("K97").Value = 549.5
$ws.Range("M97").Value = -53.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3333
$ws.Range("I82").Value = 2499.5
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 2499.5
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -2138.5
$ws.Range("N82").Value = -5722

$ws.Range("H85").Value = 3333
$ws.Range("I85").Value = 2499.5
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 2499.5
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -1251.5
$ws.Range("N85").Value = -7496

$ws.Range("H93").Value = 2672.5454
$ws.Range("I93").Value = 2279.8
$ws.Range("J93").Value = 2999.8333
$ws.Range("K93").Value = 2279.8
$ws.Range("L93").Value = 2999.8333
$ws.Range("M93").Value = -1031.8
$ws.Range("N93").Value = -5495.8333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7343.1665
$ws.Range("J81").Value = 7343.1665
$ws.Range("L81").Value = 14686.333
$ws.Range("N81").Value = -16808.333

$ws.Range("H84").Value = 7343.1665
$ws.Range("J84").Value = 7343.1665
$ws.Range("L84").Value = 73431.66500000001
$ws.Range("N84").Value = -84039.66500000001

$ws.Range("H126").Value = 1249.8334
$ws.Range("I126").Value = 1166.3334
$ws.Range("J126").Value = 1333.3334
$ws.Range("K126").Value = 3499.0002
$ws.Range("L126").Value = 4000.0002
$ws.Range("M126").Value = -1029.0002
$ws.Range("N126").Value = -8940.0002
